$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: refreshed "now" snapshot value in A1 (B1/C1 unchanged)
$ws.Range("A1").Value = 43265.430810694445

# Per-row data refresh for the running log (rows 7-31):
#  - Start Date/Time (A&B) and End Date/Time (C&D) re-randomized
#  - Distance (F) re-randomized
#  - Location (E) reassigned across the shrunk location list
$rows = @(
    @{ Row = 7; A = 43101.409160868054; C = 43101.43182290509; F = 6.668033011472151; Location = "Lummen" }
    @{ Row = 8; A = 43101.9287515625; C = 43101.96088119213; F = 9.062425027218962; Location = "Schulen" }
    @{ Row = 9; A = 43106.14966076389; C = 43106.201616782404; F = 13.82755953400363; Location = "Antwerp" }
    @{ Row = 10; A = 43119.1888405787; C = 43119.21737067129; F = 9.250777324426243; Location = "Antwerp" }
    @{ Row = 11; A = 43119.54404277778; C = 43119.5649687037; F = 6.485027374596375; Location = "Brussels" }
    @{ Row = 12; A = 43119.922652534726; C = 43119.951761331016; F = 8.544304183609116; Location = "Schulen" }
    @{ Row = 13; A = 43120.84319167824; C = 43120.8652171412; F = 7.086167242639629; Location = "Antwerp" }
    @{ Row = 14; A = 43122.76418954861; C = 43122.78184001157; F = 5.849337544422425; Location = "Heusden-Zolder" }
    @{ Row = 15; A = 43133.54016383102; C = 43133.59371707176; F = 12.996725418522178; Location = "Antwerp" }
    @{ Row = 16; A = 43137.14482313657; C = 43137.194776840275; F = 11.48072402100897; Location = "Brussels" }
    @{ Row = 17; A = 43148.97855799769; C = 43149.02191447917; F = 13.745713262714593; Location = "Brussels" }
    @{ Row = 18; A = 43151.58271003472; C = 43151.615997071756; F = 8.978062942102753; Location = "Antwerp" }
    @{ Row = 19; A = 43154.95035594908; C = 43154.98872400463; F = 10.081276366171927; Location = "Antwerp" }
    @{ Row = 20; A = 43157.88583143518; C = 43157.91775273148; F = 10.42934088903244; Location = "Brussels" }
    @{ Row = 21; A = 43161.591063865744; C = 43161.61484858796; F = 6.249876055207235; Location = "Heusden-Zolder" }
    @{ Row = 22; A = 43169.182145162034; C = 43169.23098775463; F = 13.973749309602663; Location = "Schulen" }
    @{ Row = 23; A = 43171.23680543982; C = 43171.285335532404; F = 13.479708832720782; Location = "Antwerp" }
    @{ Row = 24; A = 43172.49448210648; C = 43172.520604791665; F = 7.9199883455004825; Location = "Heusden-Zolder" }
    @{ Row = 25; A = 43180.63088365741; C = 43180.67786282407; F = 11.919680856654605; Location = "Brussels" }
    @{ Row = 26; A = 43180.848199375; C = 43180.910491041664; F = 13.632065683365433; Location = "Antwerp" }
    @{ Row = 27; A = 43184.41087211805; C = 43184.433198506944; F = 5.254952817949911; Location = "Heusden-Zolder" }
    @{ Row = 28; A = 43192.854273668985; C = 43192.89941255787; F = 12.336936916024031; Location = "Heusden-Zolder" }
    @{ Row = 29; A = 43194.31599721065; C = 43194.35825415509; F = 13.334186775526472; Location = "Antwerp" }
    @{ Row = 30; A = 43209.2122253125; C = 43209.249308645834; F = 12.46160960707984; Location = "Schulen" }
    @{ Row = 31; A = 43220.35470402778; C = 43220.37970402778; F = 7.449201245258642; Location = "Lummen" }
)

foreach ($row in $rows) {
    $r = $row.Row
    $ws.Cells.Item($r, 1).Value = $row.A   # A: Start Date
    $ws.Cells.Item($r, 2).Value = $row.A   # B: Start Time (mirrors A)
    $ws.Cells.Item($r, 3).Value = $row.C   # C: End Date
    $ws.Cells.Item($r, 4).Value = $row.C   # D: End Time (mirrors C)
    $ws.Cells.Item($r, 5).Value = $row.Location   # E: Location
    $ws.Cells.Item($r, 6).Value = $row.F   # F: Distance
}

# Turn on the table AutoFilter dropdowns (previously absent)
$lo = $ws.ListObjects("RUNNINGRESULTS")
$lo.ShowAutoFilter = $true
